$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.152.26"
$ws.Range("E2").Value = "  -3.55%  "
$ws.Range("D3").Value = "2.979.54"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "558.45"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").Value = "134.36"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").Value = "2.972.37"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").Value = "4.87"
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "33.15"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.463.72"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "6.90"
$ws.Range("E17").Value = "  +10.77%  "
$ws.Range("D18").Value = "2.971.40"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "58.030.70"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "421.79"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "0.689"
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("D23").Value = "7.02"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "13.11"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").Value = "79.79"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "2.51"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "7.65"
$ws.Range("E29").Value = "  +5.05%  "
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").Value = "  +6.79%  "
$ws.Range("D31").Value = "25.37"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "0.100"
$ws.Range("E33").Value = "  +7.19%  "
$ws.Range("D34").Value = "2.16"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "5.69"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "0.949"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "0.0₃0698"
$ws.Range("E37").Value = "  +5.63%  "
$ws.Range("D38").Value = "48.68"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").Value = "8.58"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("D40").Value = "2.60"
$ws.Range("E40").Value = "  +5.48%  "
$ws.Range("D41").Value = "384.28"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "0.0352"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "2.691.81"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "0.244"
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("D47").Value = "123.16"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "23.71"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "2.03"
$ws.Range("E51").Value = "  +0.18%  "
